# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) -- row => new F value
$sheet1Updates = @{
    2  = 254
    5  = 442
    6  = 82
    7  = 571
    9  = 6858
    10 = 164
    13 = 179
    15 = 1113
    16 = 16296
    17 = 1608
    19 = 336
    22 = 11429
    23 = 13
    24 = 1065
    25 = 4499
    26 = 355
    31 = 143
}

# Sheet "全部类型" (4th sheet) -- row => new F value
$sheet4Updates = @{
    2  = 254
    5  = 442
    6  = 82
    7  = 571
    10 = 6858
    11 = 164
    14 = 179
    17 = 1113
    18 = 16296
    19 = 1608
    21 = 336
    26 = 11429
    27 = 13
    28 = 1065
    29 = 4499
    30 = 355
    35 = 143
}

$ws1 = $wb.Worksheets.Item(1)
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item(4)
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
